$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5955
$ws.Range("E2").Value = -672
$ws.Range("F2").Value = -700
$ws.Range("G2").Value = -2236
$ws.Range("H2").Value = -2826
$ws.Range("I2").Value = -2735
$ws.Range("J2").Value = -91
$ws.Range("K2").Value = 22171
$ws.Range("L2").Value = 18434
$ws.Range("M2").Value = 3737
$ws.Range("N2").Value = 4441
$ws.Range("O2").Value = -704
$ws.Range("P2").Value = 400
$ws.Range("Q2").Value = -69
$ws.Range("R2").Value = 169
$ws.Range("S2").Value = -210
$ws.Range("T2").Value = 18
$ws.Range("U2").Value = -87
$ws.Range("V2").Value = 11392
$ws.Range("W2").Value = -11.29
$ws.Range("X2").Value = -47.46
$ws.Range("Y2").Value = -47.32
$ws.Range("Z2").Value = -12.27
$ws.Range("AA2").Value = 493.29
$ws.Range("AB2").Value = -619.64
$ws.Range("AC2").Value = -397968
$ws.Range("AD2").Value = -0.12
$ws.Range("AE2").Value = 646282
$ws.Range("AF2").Value = 0.07
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 687176

# Row 3
$ws.Range("D3").Value = 4468
$ws.Range("E3").Value = -726
$ws.Range("F3").Value = -633
$ws.Range("G3").Value = -5938
$ws.Range("H3").Value = -6330
$ws.Range("I3").Value = -6202
$ws.Range("J3").Value = -128
$ws.Range("K3").Value = 17861
$ws.Range("L3").Value = 20458
$ws.Range("M3").Value = -2596
$ws.Range("N3").Value = -1746
$ws.Range("O3").Value = -851
$ws.Range("P3").Value = 440
$ws.Range("Q3").Value = 149
$ws.Range("R3").Value = 964
$ws.Range("S3").Value = -1030
$ws.Range("T3").Value = 67
$ws.Range("U3").Value = 83
$ws.Range("V3").Value = 10936
$ws.Range("W3").Value = -16.26
$ws.Range("X3").Value = -141.68
$ws.Range("Y3").Value = -460.18
$ws.Range("Z3").Value = -31.62
$ws.Range("AA3").Value = -787.95
$ws.Range("AB3").Value = -1975.88
$ws.Range("AC3").Value = -844700
$ws.Range("AD3").Value = -0.05
$ws.Range("AE3").Value = -230961
$ws.Range("AF3").Value = -0.18
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 755870

# Row 4
$ws.Range("D4").Value = 3903
$ws.Range("E4").Value = -351
$ws.Range("F4").Value = -351
$ws.Range("G4").Value = 5195
$ws.Range("H4").Value = 2614
$ws.Range("I4").Value = 2653
$ws.Range("J4").Value = -38
$ws.Range("K4").Value = 4257
$ws.Range("L4").Value = 3191
$ws.Range("M4").Value = 1066
$ws.Range("N4").Value = 2191
$ws.Range("O4").Value = -1125
$ws.Range("P4").Value = 488
$ws.Range("Q4").Value = 979
$ws.Range("R4").Value = 979
$ws.Range("S4").Value = -1300
$ws.Range("T4").Value = 27
$ws.Range("U4").Value = 952
$ws.Range("V4").Value = 661
$ws.Range("W4").Value = -8.99
$ws.Range("X4").Value = 66.98
$ws.Range("Y4").Value = 1191.98
$ws.Range("Z4").Value = 23.64
$ws.Range("AA4").Value = 299.31
$ws.Range("AB4").Value = 1988.84
$ws.Range("AC4").Value = 9838
$ws.Range("AD4").Value = 0.11
$ws.Range("AE4").Value = 4490
$ws.Range("AF4").Value = 0.24
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 48791975

# Row 5
$ws.Range("D5").Value = 2804
$ws.Range("E5").Value = -194
$ws.Range("F5").Value = -194
$ws.Range("G5").Value = -351
$ws.Range("H5").Value = -404
$ws.Range("I5").Value = -330
$ws.Range("J5").Value = -74
$ws.Range("K5").Value = 3904
$ws.Range("L5").Value = 2647
$ws.Range("M5").Value = 1257
$ws.Range("N5").Value = 2361
$ws.Range("O5").Value = -1105
$ws.Range("P5").Value = 938
$ws.Range("Q5").Value = 36
$ws.Range("R5").Value = -399
$ws.Range("S5").Value = 254
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = 35
$ws.Range("V5").Value = 958
$ws.Range("W5").Value = -6.9
$ws.Range("X5").Value = -14.42
$ws.Range("Y5").Value = -14.52
$ws.Range("Z5").Value = -9.91
$ws.Range("AA5").Value = 210.62
$ws.Range("AB5").Value = 1014.63
$ws.Range("AC5").Value = -532
$ws.Range("AD5").Value = -2.62
$ws.Range("AE5").Value = 2516
$ws.Range("AF5").Value = 0.55
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 93833750

# Row 6
$ws.Range("D6").Value = 1774
$ws.Range("E6").Value = -94
$ws.Range("F6").Value = -94
$ws.Range("G6").Value = -503
$ws.Range("H6").Value = -515
$ws.Range("I6").Value = -510
$ws.Range("K6").Value = 4131
$ws.Range("L6").Value = 2228
$ws.Range("M6").Value = 1903
$ws.Range("N6").Value = 1972
$ws.Range("P6").Value = 1240
$ws.Range("Q6").Value = -178
$ws.Range("R6").Value = 258
$ws.Range("S6").Value = 236
$ws.Range("T6").Value = 124
$ws.Range("U6").Value = -302
$ws.Range("V6").Value = 844
$ws.Range("W6").Value = -5.28
$ws.Range("X6").Value = -29.02
$ws.Range("Y6").Value = -23.52
$ws.Range("Z6").Value = -12.82
$ws.Range("AA6").Value = 117.11
$ws.Range("AB6").Value = 735.48
$ws.Range("AC6").Value = -492
$ws.Range("AD6").Value = -1.76
$ws.Range("AE6").Value = 1590
$ws.Range("AF6").Value = 0.54
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 124015655

# Row 6: AI6 removed entirely
$ws.Range("AI6").ClearContents()

# Rows 7-9: clear all data columns D:AI, leaving only A,B,C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
